$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) — update 想去人数 (attendance) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 105
$ws1.Range("F6").Value = 13

# Sheet "全部类型" (All types) — same underlying data, same updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 105
$ws4.Range("F6").Value = 13
